$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 20, shifting existing rows 20-39 down to 21-40.
$ws.Rows.Item(20).Insert()

# Populate the newly inserted row 20 with the new weekly price entry.
$ws.Range("A20").Value = 11
$ws.Range("B20").Value = "Vega Monumental Concepción"
$ws.Range("C20").Value = "Bíobío"
$ws.Range("D20").Value = 44721
$ws.Range("E20").Value = 8
$ws.Range("F20").Value = 100112013
$ws.Range("G20").Value = "Alcachofa"
$ws.Range("H20").Value = "Española"
$ws.Range("I20").Value = "Primera"
$ws.Range("J20").Value = 150
$ws.Range("K20").Value = 19000
$ws.Range("L20").Value = 20000
$ws.Range("M20").Value = 19533
$ws.Range("N20").Value = "$/caja 30 unidades"
$ws.Range("O20").Value = "Provincia de Limarí"
$ws.Range("P20").Value = 651
$ws.Range("Q20").Value = 30
$ws.Range("R20").Value = "Hortaliza"
